$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

$row2 = New-Object 'object[,]' 1,7
$row2[0,0] = 'BRVM - CONSOMMATION DE BASE     (**)'
$row2[0,1] = 0
$row2[0,2] = 4
$row2[0,3] = 1132.37
$row2[0,4] = 282.33
$row2[0,5] = '🟡 Observer'
$row2[0,6] = '➖ Neutre'
$ws1.Range("A2:G2").Value = $row2

$row3 = New-Object 'object[,]' 1,7
$row3[0,0] = 'BRVM-PRINCIPAL     (**)'
$row3[0,1] = 0
$row3[0,2] = 4
$row3[0,3] = 1081.28
$row3[0,4] = 272.05
$row3[0,5] = '🟡 Observer'
$row3[0,6] = '➖ Neutre'
$ws1.Range("A3:G3").Value = $row3

$row4 = New-Object 'object[,]' 1,7
$row4[0,0] = 'BRVM - INDUSTRIELS'
$row4[0,1] = 0
$row4[0,2] = 4
$row4[0,3] = 786.69
$row4[0,4] = 208.12
$row4[0,5] = '🟡 Observer'
$row4[0,6] = '➖ Neutre'
$ws1.Range("A4:G4").Value = $row4

$row5 = New-Object 'object[,]' 1,7
$row5[0,0] = 'BRVM - CONSOMMATION DISCRETIONNAIRE'
$row5[0,1] = 0
$row5[0,2] = 4
$row5[0,3] = 767.12
$row5[0,4] = 194.74
$row5[0,5] = '🟡 Observer'
$row5[0,6] = '➖ Neutre'
$ws1.Range("A5:G5").Value = $row5

$row6 = New-Object 'object[,]' 1,7
$row6[0,0] = 'BRVM - SERVICES FINANCIERS'
$row6[0,1] = 0
$row6[0,2] = 4
$row6[0,3] = 660.75
$row6[0,4] = 165.83
$row6[0,5] = '🟡 Observer'
$row6[0,6] = '➖ Neutre'
$ws1.Range("A6:G6").Value = $row6

$row7 = New-Object 'object[,]' 1,7
$row7[0,0] = 'BRVM-PRESTIGE'
$row7[0,1] = 0
$row7[0,2] = 4
$row7[0,3] = 619.24
$row7[0,4] = 155.42
$row7[0,5] = '🟡 Observer'
$row7[0,6] = '➖ Neutre'
$ws1.Range("A7:G7").Value = $row7

$row8 = New-Object 'object[,]' 1,7
$row8[0,0] = 'BRVM – COMPOSITE TOTAL RETURN     (**)'
$row8[0,1] = 0
$row8[0,2] = 4
$row8[0,3] = 604.79
$row8[0,4] = 151.83
$row8[0,5] = '🟡 Observer'
$row8[0,6] = '➖ Neutre'
$ws1.Range("A8:G8").Value = $row8

$row9 = New-Object 'object[,]' 1,7
$row9[0,0] = 'BRVM - ENERGIE'
$row9[0,1] = 0
$row9[0,2] = 4
$row9[0,3] = 517.85
$row9[0,4] = 135.04
$row9[0,5] = '🟡 Observer'
$row9[0,6] = '➖ Neutre'
$ws1.Range("A9:G9").Value = $row9

$row10 = New-Object 'object[,]' 1,7
$row10[0,0] = 'BRVM - SERVICES PUBLICS'
$row10[0,1] = 0
$row10[0,2] = 4
$row10[0,3] = 487.74
$row10[0,4] = 124.06
$row10[0,5] = '🟡 Observer'
$row10[0,6] = '➖ Neutre'
$ws1.Range("A10:G10").Value = $row10

$row11 = New-Object 'object[,]' 1,7
$row11[0,0] = 'BRVM - TELECOMMUNICATIONS'
$row11[0,1] = 0
$row11[0,2] = 4
$row11[0,3] = 401.83
$row11[0,4] = 100.45
$row11[0,5] = '🟡 Observer'
$row11[0,6] = '➖ Neutre'
$ws1.Range("A11:G11").Value = $row11

$row12 = New-Object 'object[,]' 1,7
$row12[0,0] = 'SAFCA CI (SAFC)'
$row12[0,1] = 2
$row12[0,2] = 0
$row12[0,3] = 14.92
$row12[0,4] = 7.43
$row12[0,5] = '🟡 Observer'
$row12[0,6] = '➖ Neutre'
$ws1.Range("A12:G12").Value = $row12

$row13 = New-Object 'object[,]' 1,7
$row13[0,0] = 'SERVAIR ABIDJAN CI (ABJC)'
$row13[0,1] = 1
$row13[0,2] = 0
$row13[0,3] = 14.81
$row13[0,4] = 7.46
$row13[0,5] = '🟡 Observer'
$row13[0,6] = '➖ Neutre'
$ws1.Range("A13:G13").Value = $row13

$row14 = New-Object 'object[,]' 1,7
$row14[0,0] = 'SUCRIVOIRE (SCRC)'
$row14[0,1] = 2
$row14[0,2] = 1
$row14[0,3] = 12.68
$row14[0,4] = -2.15
$row14[0,5] = '🟡 Observer'
$row14[0,6] = '👀 À surveiller'
$ws1.Range("A14:G14").Value = $row14

$row15 = New-Object 'object[,]' 1,7
$row15[0,0] = 'ECOBANK TRANS. INCORP. TG (ETIT)'
$row15[0,1] = 2
$row15[0,2] = 1
$row15[0,3] = 10.41
$row15[0,4] = -2.94
$row15[0,5] = '🟡 Observer'
$row15[0,6] = '👀 À surveiller'
$ws1.Range("A15:G15").Value = $row15

$row16 = New-Object 'object[,]' 1,7
$row16[0,0] = 'ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)'
$row16[0,1] = 1
$row16[0,2] = 0
$row16[0,3] = 7.45
$row16[0,4] = 7.45
$row16[0,5] = '🟡 Observer'
$row16[0,6] = '➖ Neutre'
$ws1.Range("A16:G16").Value = $row16

$row17 = New-Object 'object[,]' 1,7
$row17[0,0] = 'EVIOSYS PACKAGING SIEM CI (SEMC)'
$row17[0,1] = 2
$row17[0,2] = 1
$row17[0,3] = 7.44
$row17[0,4] = 7.33
$row17[0,5] = '🟡 Observer'
$row17[0,6] = '👀 À surveiller'
$ws1.Range("A17:G17").Value = $row17

$row18 = New-Object 'object[,]' 1,7
$row18[0,0] = 'NEI-CEDA CI (NEIC)'
$row18[0,1] = 1
$row18[0,2] = 0
$row18[0,3] = 6.92
$row18[0,4] = 6.92
$row18[0,5] = '🟡 Observer'
$row18[0,6] = '➖ Neutre'
$ws1.Range("A18:G18").Value = $row18

$row19 = New-Object 'object[,]' 1,7
$row19[0,0] = 'TOTALENERGIES MARKETING CI (TTLC)'
$row19[0,1] = 1
$row19[0,2] = 0
$row19[0,3] = 6.75
$row19[0,4] = 6.75
$row19[0,5] = '🟡 Observer'
$row19[0,6] = '➖ Neutre'
$ws1.Range("A19:G19").Value = $row19

$row20 = New-Object 'object[,]' 1,7
$row20[0,0] = 'ONATEL BF (ONTBF)'
$row20[0,1] = 1
$row20[0,2] = 0
$row20[0,3] = 6.3
$row20[0,4] = 6.3
$row20[0,5] = '🟡 Observer'
$row20[0,6] = '➖ Neutre'
$ws1.Range("A20:G20").Value = $row20

$row21 = New-Object 'object[,]' 1,7
$row21[0,0] = 'CORIS BANK INTERNATIONAL (CBIBF)'
$row21[0,1] = 1
$row21[0,2] = 1
$row21[0,3] = 5.95
$row21[0,4] = -1.55
$row21[0,5] = '🟡 Observer'
$row21[0,6] = '👀 À surveiller'
$ws1.Range("A21:G21").Value = $row21

$row22 = New-Object 'object[,]' 1,7
$row22[0,0] = 'SETAO CI (STAC)'
$row22[0,1] = 1
$row22[0,2] = 1
$row22[0,3] = 4.36
$row22[0,4] = 7.35
$row22[0,5] = '🟡 Observer'
$row22[0,6] = '👀 À surveiller'
$ws1.Range("A22:G22").Value = $row22

$row23 = New-Object 'object[,]' 1,7
$row23[0,0] = 'SITAB CI (STBC)'
$row23[0,1] = 1
$row23[0,2] = 1
$row23[0,3] = 3.39
$row23[0,4] = -1.9
$row23[0,5] = '🟡 Observer'
$row23[0,6] = '👀 À surveiller'
$ws1.Range("A23:G23").Value = $row23

$row24 = New-Object 'object[,]' 1,7
$row24[0,0] = 'UNILEVER CI (UNLC)'
$row24[0,1] = 2
$row24[0,2] = 1
$row24[0,3] = 0.46
$row24[0,4] = -6.42
$row24[0,5] = '🟡 Observer'
$row24[0,6] = '👀 À surveiller'
$ws1.Range("A24:G24").Value = $row24

$row25 = New-Object 'object[,]' 1,7
$row25[0,0] = 'SMB CI (SMBC)'
$row25[0,1] = 2
$row25[0,2] = 1
$row25[0,3] = 0.3
$row25[0,4] = 6.66
$row25[0,5] = '🟡 Observer'
$row25[0,6] = '👀 À surveiller'
$ws1.Range("A25:G25").Value = $row25

$row27 = New-Object 'object[,]' 1,7
$row27[0,0] = 'SICABLE CI (CABC)'
$row27[0,1] = 1
$row27[0,2] = 2
$row27[0,3] = -2.36
$row27[0,4] = 7.46
$row27[0,5] = '🟡 Observer'
$row27[0,6] = '👀 À surveiller'
$ws1.Range("A27:G27").Value = $row27

$row28 = New-Object 'object[,]' 1,7
$row28[0,0] = 'SOLIBRA CI (SLBC)'
$row28[0,1] = 0
$row28[0,2] = 1
$row28[0,3] = -2.75
$row28[0,4] = -2.75
$row28[0,5] = '🟡 Observer'
$row28[0,6] = '➖ Neutre'
$ws1.Range("A28:G28").Value = $row28

$row29 = New-Object 'object[,]' 1,7
$row29[0,0] = 'SOCIETE IVOIRIENNE DE BANQUE  (SIBC)'
$row29[0,1] = 0
$row29[0,2] = 1
$row29[0,3] = -4.1
$row29[0,4] = -4.1
$row29[0,5] = '🟡 Observer'
$row29[0,6] = '➖ Neutre'
$ws1.Range("A29:G29").Value = $row29

$row31 = New-Object 'object[,]' 1,7
$row31[0,0] = 'ORANGE COTE D''IVOIRE (ORAC)'
$row31[0,1] = 0
$row31[0,2] = 1
$row31[0,3] = -6.93
$row31[0,4] = -6.93
$row31[0,5] = '🟡 Observer'
$row31[0,6] = '➖ Neutre'
$ws1.Range("A31:G31").Value = $row31

$row32 = New-Object 'object[,]' 1,7
$row32[0,0] = 'ECOBANK COTE D''''IVOIRE (ECOC)'
$row32[0,1] = 0
$row32[0,2] = 2
$row32[0,3] = -9.74
$row32[0,4] = -4.41
$row32[0,5] = '🟡 Observer'
$row32[0,6] = '➖ Neutre'
$ws1.Range("A32:G32").Value = $row32

$row33 = New-Object 'object[,]' 1,7
$row33[0,0] = 'SICOR CI (SICC)'
$row33[0,1] = 0
$row33[0,2] = 2
$row33[0,3] = -14.14
$row33[0,4] = -6.72
$row33[0,5] = '🟡 Observer'
$row33[0,6] = '➖ Neutre'
$ws1.Range("A33:G33").Value = $row33

$ws2.Range("B2").Value = 21438.12
$ws2.Range("B3").Value = 18705.17
$ws2.Range("B4").Value = 7633.38
$ws2.Range("B5").Value = 7146.85
$ws2.Range("B6").Value = 4845.3
$ws2.Range("B7").Value = 4115.64
$ws2.Range("B8").Value = 3881.6
$ws2.Range("B9").Value = 2671.17
$ws2.Range("B10").Value = 2325.79
$ws2.Range("B11").Value = 1514.61
